$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Column A (Category) values to the new, consolidated category labels.
# The order below matches the grouping the values were authored in, so that
# newly introduced shared strings come out in the same order as the source
# workbook: Cup Dhoop, HEM Camphor, Pooja Oil, Ghee Diya, Pooja Samagri.
$ws.Range("A12").Value = " Cup Dhoop "
$ws.Range("A13").Value = " Cup Dhoop "
$ws.Range("A14").Value = " Cup Dhoop "

$ws.Range("A5").Value = "HEM  Camphor "
$ws.Range("A6").Value = "HEM  Camphor "

$ws.Range("A3").Value = " Pooja Oil "
$ws.Range("A4").Value = " Pooja Oil "

$ws.Range("A7").Value = "Ghee Diya"
$ws.Range("A8").Value = "Ghee Diya"
$ws.Range("A9").Value = "Ghee Diya"
$ws.Range("A10").Value = "Ghee Diya"

$ws.Range("A2").Value = "Pooja Samagri "
$ws.Range("A11").Value = "Pooja Samagri "

# A11 previously used a slightly different cell style (no explicit font
# applied). Align it with the rest of column A by copying the format from
# the cell directly above it.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Update the active selection to match the new cursor position.
$ws.Range("A7").Select()
